$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Helper: replace the single paragraph inside a table cell with an exact
# sequence of run/proofErr XML, giving byte-for-byte control over the
# resulting WordprocessingML (run splits, xml:space, proofErr markers, ...).
function Set-CellParagraphXml($row, $col, $innerXml) {
  $cell = $t.Cell($row, $col)
  $range = $cell.Range
  $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'
  [void]$range.InsertXML($pkg)
}

# 1) Shade the title cell (row 1) with the accent1/99% theme fill (closest
#    reachable approximation via the Shading object model: clear pattern,
#    automatic foreground, RGB background matching accent1 @ 99% tint).
$titleCell = $t.Cell(1, 1)
$titleCell.Shading.Texture = 0            # wdTextureNone -> w:val="clear"
$titleCell.Shading.ForegroundPatternColor = -16777216   # wdColorAutomatic -> w:color="auto"
$titleCell.Shading.BackgroundPatternColor = 14396046    # RGB(142,170,219) = 0x8EAADB -> w:fill

# 2) Fix "Ing.Jenny" -> "Ing. Jenny" (and drop the spell-check markers) in
#    the Tutora row. Because that row's two cells use gridSpan="2", the
#    second real <w:tc> is reached via column index 3.
Set-CellParagraphXml 4 3 '<w:r><w:t>Ing. Jenny</w:t></w:r><w:r><w:t xml:space="preserve"> Ruiz</w:t></w:r>'

# 3) Fill in the previously-empty test case rows.
Set-CellParagraphXml 8 1 '<w:r><w:t>REQ001</w:t></w:r>'
Set-CellParagraphXml 8 2 '<w:r><w:t>Ingresar nuevo usuario</w:t></w:r>'
Set-CellParagraphXml 8 3 '<w:r><w:t>Pendiente ventana de registro</w:t></w:r>'
Set-CellParagraphXml 8 4 '<w:r><w:t xml:space="preserve">Diseño de </w:t></w:r><w:r><w:t>página</w:t></w:r><w:r><w:t xml:space="preserve"> de CSS</w:t></w:r>'

Set-CellParagraphXml 9 1 '<w:r><w:t>REQ0</w:t></w:r><w:r><w:t>02</w:t></w:r>'
Set-CellParagraphXml 9 2 '<w:r><w:t>Inicio de sesión y contraseña</w:t></w:r>'
Set-CellParagraphXml 9 3 '<w:r><w:t>Validación de datos para ingreso</w:t></w:r>'
Set-CellParagraphXml 9 4 '<w:r><w:t xml:space="preserve">Investigar CSS para generar los </w:t></w:r><w:r><w:t xml:space="preserve">datos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>y</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ingresar</w:t></w:r>'

Set-CellParagraphXml 10 1 '<w:r><w:t>REQ00</w:t></w:r><w:r><w:t>3</w:t></w:r>'
Set-CellParagraphXml 10 2 '<w:r><w:t>Menú e iconos</w:t></w:r>'
Set-CellParagraphXml 10 3 '<w:r><w:t>Los iconos no direccionan</w:t></w:r>'
Set-CellParagraphXml 10 4 '<w:r><w:t>Investigando CSS para direccionar</w:t></w:r>'

Set-CellParagraphXml 11 1 '<w:r><w:t>REQ00</w:t></w:r><w:r><w:t>4</w:t></w:r>'
Set-CellParagraphXml 11 2 '<w:r><w:t xml:space="preserve">Diseño </w:t></w:r><w:r><w:t>página</w:t></w:r><w:r><w:t xml:space="preserve"> amigable</w:t></w:r>'
Set-CellParagraphXml 11 3 '<w:r><w:t>Cambiar tonalidades y texto</w:t></w:r>'
Set-CellParagraphXml 11 4 '<w:r><w:t xml:space="preserve">Se esta solicitando al cliente </w:t></w:r><w:r><w:t>los matices</w:t></w:r><w:r><w:t xml:space="preserve"> de colores que requiere</w:t></w:r>'

# 4) Append a brand-new row (REQ005) at the end of the table.
$t.Rows.Add() | Out-Null
$newRowIndex = $t.Rows.Count
Set-CellParagraphXml $newRowIndex 1 '<w:r><w:t>REQ005</w:t></w:r>'
Set-CellParagraphXml $newRowIndex 2 '<w:r><w:t>Productos</w:t></w:r>'
Set-CellParagraphXml $newRowIndex 3 '<w:r><w:t>Falta añadir contenido y enlaces</w:t></w:r>'
Set-CellParagraphXml $newRowIndex 4 '<w:r><w:t>Se debe generar una nueva ventana</w:t></w:r>'

Write-Host "Edit complete"
